# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml (the "Integral" / Red Violet theme used by the
#   slide master) is replaced with the "Office Theme" color scheme that
#   used to live in ppt/theme/theme2.xml (the notes-master theme).
#
# This headless COM host only exposes theme editing through the
# SlideMaster's ColorScheme (RGBColor.RGB) - there is no object model
# path that reaches the notes-master's theme part independently, so we
# apply the reachable half of the swap: push the "Office Theme" palette
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink, in PowerPoint's standard
# 1-12 ColorScheme.Colors ordering) onto the slide master's color scheme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

# RGB values below are packed as 0xBBGGRR, matching PowerPoint's
# RGBColor.RGB convention (same values COM callers use for
# Colors(i).RGB = ...).
$scheme.Colors(1).RGB  = 0        # dk1      000000
$scheme.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388  # dk2      44546A
$scheme.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501  # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407    # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308 # accent5  4472C4
$scheme.Colors(10).RGB = 4697456  # accent6  70AD47
$scheme.Colors(11).RGB = 12673797 # hlink    0563C1
$scheme.Colors(12).RGB = 7491477  # folHlink 954F72
